$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Bad Drivers" table (weekly refreshed numbers) ---
$ws.Range("C3").Value = 1804
$ws.Range("D3").Value = 61.8

$ws.Range("C4").Value = 4083
$ws.Range("D4").Value = 92.2

$ws.Range("D5").Value = 96.8

$ws.Range("C6").Value = 6227

# --- "Good Drivers" table: weekly refresh shifts every row down by one,
# introduces a newly-qualified driver at the top, drops the oldest one
# off the bottom, and picks up newly-learned "Driver Vintage" dates for
# a couple of rows that previously had none. Write the full new table
# contents row by row (A/B/D always; E only when a vintage date is known).

function Set-VintageDate($cellRange, $copyFromRange, $dateText) {
    if ($dateText -eq $null) {
        $cellRange.ClearContents()
    } else {
        $cellRange.Value = "'" + $dateText
        $copyFromRange.Copy()
        $cellRange.PasteSpecial(-4122)
    }
}

# Row 14
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B14").Value = 445055
$ws.Range("D14").Value = 99.90000000000001
Set-VintageDate $ws.Range("E14") $ws.Range("D14") "2024-11-10"

# Row 15
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1"
$ws.Range("B15").Value = 10661
$ws.Range("D15").Value = 100
Set-VintageDate $ws.Range("E15") $ws.Range("D15") "2022-08-29"

# Row 16
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3"
$ws.Range("B16").Value = 14239
$ws.Range("D16").Value = 100
Set-VintageDate $ws.Range("E16") $ws.Range("D16") "2022-05-23"

# Row 17
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1"
$ws.Range("B17").Value = 265400
$ws.Range("D17").Value = 99.90000000000001
Set-VintageDate $ws.Range("E17") $ws.Range("D17") "2022-05-01"

# Row 18
$ws.Range("A18").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B18").Value = 77849
$ws.Range("D18").Value = 99.90000000000001
Set-VintageDate $ws.Range("E18") $ws.Range("D18") "2021-08-18"

# Row 19
$ws.Range("A19").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B19").Value = 34244
$ws.Range("D19").Value = 100
Set-VintageDate $ws.Range("E19") $ws.Range("D19") "2021-04-27"

# Row 20
$ws.Range("A20").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B20").Value = 59673
$ws.Range("D20").Value = 100
Set-VintageDate $ws.Range("E20") $ws.Range("D20") "2020-08-05"

# Row 21
$ws.Range("A21").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B21").Value = 113652
$ws.Range("D21").Value = 100
Set-VintageDate $ws.Range("E21") $ws.Range("D21") "2020-01-06"

# Row 22
$ws.Range("A22").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B22").Value = 56018
$ws.Range("D22").Value = 100
Set-VintageDate $ws.Range("E22") $ws.Range("D22") $null
